$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 81 - this shifts the existing rows 81..112 down to 82..113
$ws.Rows.Item(81).Insert()

# Copy the date cell's number format from the row below (now row 82) so the
# new D81 date renders the same way as the rest of the column.
$ws.Cells.Item(81, 4).NumberFormat = $ws.Cells.Item(82, 4).NumberFormat

# Populate the new weekly record in row 81. The non-varying (descriptive)
# columns are identical to every other row in this data set.
$ws.Cells.Item(81, 1).Value = 6
$ws.Cells.Item(81, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(81, 3).Value = "Metropolitana"
$ws.Cells.Item(81, 4).Value = 44900
$ws.Cells.Item(81, 5).Value = 13
$ws.Cells.Item(81, 6).Value = 100114007
$ws.Cells.Item(81, 7).Value = "Jengibre"
$ws.Cells.Item(81, 8).Value = "Sin especificar"
$ws.Cells.Item(81, 9).Value = "Primera"
$ws.Cells.Item(81, 10).Value = 400
$ws.Cells.Item(81, 11).Value = 10000
$ws.Cells.Item(81, 12).Value = 11000
$ws.Cells.Item(81, 13).Value = 10425
$ws.Cells.Item(81, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(81, 15).Value = "Perú"
$ws.Cells.Item(81, 16).Value = 802
$ws.Cells.Item(81, 17).Value = 13
$ws.Cells.Item(81, 18).Value = "Hortaliza"
